# vtp, hsrp, ether channel, raah
#
# Applies:
#  1) Sheet "IP cimek" (sheet1): finish out the 192.168.2.x /30 block in row 57
#     (D:H) and append the new 192.168.3.0/27 subnetting breakdown as rows
#     60-64.
#  2) Sheet "Hatarido k" (sheet2): mark vlanok / vtp / router-on-a-stick /
#     hsrp as done (orange highlight on C6:C9), matching the commit message.
#  3) View-state touch-up: sheet1 becomes the active/selected tab, zoomed to
#     70%, with the cursor left on J63; sheet2's cursor is left on C11.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1: "IP cimek" -----------------------------------------------

# Finish row 57 (192.168.2.160/30 block): first usable, last usable,
# broadcast, mask, CIDR.
$ws1.Cells.Item(57, 4).Value = "192.168.2.161"
$ws1.Cells.Item(57, 5).Value = "192.168.2.162"
$ws1.Cells.Item(57, 6).Value = "192.168.2.163"
$ws1.Cells.Item(57, 7).Value = "255.255.255.252"
$ws1.Cells.Item(57, 8).Value = "/30"

# New 192.168.3.0/27 subnet breakdown (data / voice / management / wireless).
$ws1.Rows.Item(60).RowHeight = 20.1
$ws1.Cells.Item(60, 1).Value = "data"
$ws1.Cells.Item(60, 2).Value = 16
$ws1.Cells.Item(60, 3).Value = "192.168.3.0"
$ws1.Cells.Item(60, 4).Value = 1
$ws1.Cells.Item(60, 5).Value = 14
$ws1.Cells.Item(60, 6).Value = "192.168.3.15"
$ws1.Cells.Item(60, 7).Value = 240

$ws1.Rows.Item(61).RowHeight = 20.1
$ws1.Cells.Item(61, 1).Value = "voice"
$ws1.Cells.Item(61, 2).Value = 8
$ws1.Cells.Item(61, 3).Value = "192.168.3.16"
$ws1.Cells.Item(61, 4).Value = 17
$ws1.Cells.Item(61, 5).Value = 22
$ws1.Cells.Item(61, 6).Value = "192.168.3.23"
$ws1.Cells.Item(61, 7).Value = 248

$ws1.Rows.Item(62).RowHeight = 20.1
$ws1.Cells.Item(62, 1).Value = "management"
$ws1.Cells.Item(62, 2).Value = 4
$ws1.Cells.Item(62, 3).Value = "192.168.3.24"
$ws1.Cells.Item(62, 4).Value = 25
$ws1.Cells.Item(62, 5).Value = 26
$ws1.Cells.Item(62, 6).Value = "192.168.3.27"
$ws1.Cells.Item(62, 7).Value = 252

$ws1.Rows.Item(63).RowHeight = 20.1
$ws1.Cells.Item(63, 1).Value = "wireless"
$ws1.Cells.Item(63, 2).Value = 4
$ws1.Cells.Item(63, 3).Value = "192.168.3.28"
$ws1.Cells.Item(63, 4).Value = 29
$ws1.Cells.Item(63, 5).Value = 30
$ws1.Cells.Item(63, 6).Value = "192.168.3.31"
$ws1.Cells.Item(63, 7).Value = 252

$ws1.Rows.Item(64).RowHeight = 20.1
$ws1.Cells.Item(64, 3).Value = "192.168.3.32"

# --- Sheet2: "Hatarido k" ----------------------------------------------

# vlanok / vtp / router on a stick / hsrp done -> highlight orange.
$ws2.Range("C6:C9").Interior.Color = 49407

# --- View state ----------------------------------------------------------

# Leave the cursor on C11 in sheet2 before switching away from it, so the
# saved selection matches (last-touched sheet keeps focus otherwise).
$ws2.Range("C11").Select() | Out-Null

$ws1.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 70
$ws1.Range("J63").Select() | Out-Null
